# Update workbook: the oldest pending case (row 2, Caso -18 / GARCIA TEODORO)
# was resolved/removed from the tracking sheet. Delete that entire row so
# every subsequent record shifts up by one, matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

$ws.Rows.Item(2).Delete()
